$wb = $excel.ActiveWorkbook

# --- Sheet "Giftcard Payments" (xl/worksheets/sheet22.xml) ---
$ws = $wb.Worksheets.Item("Giftcard Payments")

# Insert a new column before column I, and a new row before row 5.
$ws.Columns("I").Insert()
$ws.Rows(5).Insert()

# Fill in the new column / row data (Pro user -> Osprey changes).
$ws.Range("I1").Value = "Osprey"
$ws.Range("I3").Value = "Osprey Gift Cards"
$ws.Range("A5").Value = "ProShippingMethod"
$ws.Range("AE5").Formula = "' Sale"

# Text correction elsewhere on the same sheet (SKIMMER 28 -> Skimmer 28).
$ws.Range("Z14").Value = "Skimmer 28"

# Update the view state for this sheet.
$ws.Activate()
$ws.Range("AC7").Select()
$excel.ActiveWindow.ScrollColumn = 5

# --- Sheet "PDP" (xl/worksheets/sheet19.xml) ---
$wsPdp = $wb.Worksheets.Item("PDP")
$wsPdp.Activate()
$wsPdp.Range("AC7").Select()
$excel.ActiveWindow.ScrollColumn = 4

# Restore the active sheet back to Giftcard Payments (matches the saved workbook view).
$ws.Activate()
